# Update column F (dSF) values for the manaea_sean.xlsx workbook.
# This reflects a "repull data, push all data, mean calculation" update
# where the dSF column receives refreshed values while dS0 stays intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = 3
    4  = 7
    9  = -2
    10 = 5
    15 = 2
    16 = -9
    17 = 2
    18 = 1
    19 = -4
    20 = -4
    22 = -2
    24 = -4
    26 = 6
    30 = 0
    31 = -6
    32 = -2
    33 = -2
    35 = -4
    36 = -1
    37 = 3
    38 = 4
    39 = 4
    40 = 3
    41 = 1
    42 = -1
    43 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
